# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped values, as produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain TEXT (the Price column holds
# strings like "68.231.13" / "27.39" which would otherwise be auto-converted
# to numbers by a normal .Value assignment). Writing it as a formula that
# evaluates to the literal string, then collapsing it to a static value via
# Copy/PasteSpecial(xlPasteValues), keeps the cell text-typed without touching
# its number format / style.
function Set-TextValue {
    param([string]$addr, [string]$value)
    $escaped = $value -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# --- Price column (D) ---
    Set-TextValue "D2" '68.231.13'
    Set-TextValue "D3" '3.272.72'
    Set-TextValue "D5" '584.18'
    Set-TextValue "D6" '184.83'
    Set-TextValue "D8" '0.602'
    Set-TextValue "D12" '3.841.80'
    Set-TextValue "D14" '68.189.02'
    Set-TextValue "D15" '27.39'
    Set-TextValue "D17" '3.276.47'
    Set-TextValue "D19" '13.29'
    Set-TextValue "D20" '416.56'
    Set-TextValue "D23" '71.09'
    Set-TextValue "D24" '0.509'
    Set-TextValue "D30" '22.68'
    Set-TextValue "D31" '5.44'
    Set-TextValue "D32" '1.25'
    Set-TextValue "D34" '164.26'
    Set-TextValue "D36" '1.89'
    Set-TextValue "D37" '26.64'
    Set-TextValue "D39" '4.45'
    Set-TextValue "D40" '6.27'
    Set-TextValue "D41" '2.639.83'
    Set-TextValue "D44" '334.58'
    Set-TextValue "D45" '24.21'
    Set-TextValue "D47" '0.988'
    Set-TextValue "D50" '30.71'

# --- Volume(1h) column (E) ---
    $ws.Range("E2").Value = '  +0.21%  '
    $ws.Range("E3").Value = '  +0.53%  '
    $ws.Range("E4").Value = '  -0.03%  '
    $ws.Range("E5").Value = '  +0.30%  '
    $ws.Range("E6").Value = '  +0.12%  '
    $ws.Range("E8").Value = '  +0.66%  '
    $ws.Range("E9").Value = '  -1.62%  '
    $ws.Range("E10").Value = '  -0.24%  '
    $ws.Range("E11").Value = '  -2.69%  '
    $ws.Range("E12").Value = '  +0.50%  '
    $ws.Range("E13").Value = '  +1.02%  '
    $ws.Range("E14").Value = '  -0.01%  '
    $ws.Range("E15").Value = '  -2.65%  '
    $ws.Range("E16").Value = '  -1.55%  '
    $ws.Range("E17").Value = '  -0.46%  '
    $ws.Range("E18").Value = '  -2.01%  '
    $ws.Range("E19").Value = '  -2.15%  '
    $ws.Range("E20").Value = '  +6.09%  '
    $ws.Range("E21").Value = '  -2.12%  '
    $ws.Range("E22").Value = '  +0.08%  '
    $ws.Range("E23").Value = '  -0.53%  '
    $ws.Range("E24").Value = '  -1.91%  '
    $ws.Range("E25").Value = '  -1.70%  '
    $ws.Range("E26").Value = '  -0.71%  '
    $ws.Range("E27").Value = '  -3.66%  '
    $ws.Range("E28").Value = '  -0.24%  '
    $ws.Range("E29").Value = '  -1.65%  '
    $ws.Range("E30").Value = '  -0.91%  '
    $ws.Range("E31").Value = '  -4.22%  '
    $ws.Range("E32").Value = '  -2.66%  '
    $ws.Range("E33").Value = '  -4.25%  '
    $ws.Range("E34").Value = '  +1.05%  '
    $ws.Range("E35").Value = '  -3.78%  '
    $ws.Range("E36").Value = '  -3.50%  '
    $ws.Range("E37").Value = '  -0.59%  '
    $ws.Range("E38").Value = '  -3.73%  '
    $ws.Range("E39").Value = '  -3.10%  '
    $ws.Range("E40").Value = '  -3.89%  '
    $ws.Range("E41").Value = '  -0.56%  '
    $ws.Range("E42").Value = '  -1.98%  '
    $ws.Range("E43").Value = '  -3.24%  '
    $ws.Range("E44").Value = '  -1.25%  '
    $ws.Range("E45").Value = '  -4.77%  '
    $ws.Range("E46").Value = '  -2.91%  '
    $ws.Range("E47").Value = '  +0.11%  '
    $ws.Range("E48").Value = '  -2.12%  '
    $ws.Range("E49").Value = '  -1.36%  '
    $ws.Range("E50").Value = '  -2.36%  '
    $ws.Range("E51").Value = '  +0.00%  '

$excel.CutCopyMode = $false

